# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted into the Mango sheet at row 382,
# pushing the existing rows 382-424 down to 383-425 (dimension grows from
# A1:T424 to A1:T425). The new row carries a fresh observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 382; existing rows 382..424 shift to 383..425.
# The inserted row inherits formatting (incl. the date style) from the row above.
$ws.Rows(382).Insert()

$ws.Range("A382").Value = 4
$ws.Range("B382").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C382").Value = "Los Lagos"
$ws.Range("D382").Value = 45142
$ws.Range("E382").Value = 10
$ws.Range("F382").Value = "Fruta"
$ws.Range("G382").Value = 100108
$ws.Range("H382").Value = "Tropicales y subtropicales"
$ws.Range("I382").Value = 100108002
$ws.Range("J382").Value = "Mango"
$ws.Range("K382").Value = "Sin especificar"
$ws.Range("L382").Value = "Primera"
$ws.Range("M382").Value = 100
$ws.Range("N382").Value = 10000
$ws.Range("O382").Value = 10000
$ws.Range("P382").Value = 10000
$ws.Range("Q382").Value = "`$/bandeja 4 kilos"
$ws.Range("R382").Value = "Perú"
$ws.Range("S382").Value = 2500
$ws.Range("T382").Value = 4
